# PNW_generators.xlsx edit: "generatos file, hydro data"
#
# - Gen_2011 (sheet2) row 96: G/H/I become flat 10s (was literature-derived
#   heat-rate curve), and O96 picks up the same red "estimate" font used on
#   G96:K96.
# - Five new hydro-import rows (97-101) are appended: P3I, P8I, P14I, P65I,
#   P66I, each typed "imports" / zone "PNW", with zero capacity/cost and
#   mustrun flags (L=M=1) like the other PNW_HYDRO placeholder rows.
# - Selection moves down to reflect the newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gen_2011")
$ws.Activate()

# --- row 96: heat-rate segments flattened to 10, O96 re-flagged red ---
$ws.Range("G96").Value = 10
$ws.Range("H96").Value = 10
$ws.Range("I96").Value = 10
$ws.Range("O96").Font.Color = $ws.Range("G96").Font.Color

# --- new rows 97-101: PNW hydro import placeholders ---
# Names first (so the shared-string table gets P3I/P8I/P14I/P65I/P66I in
# row order), then the repeated "imports" class string.
$ws.Range("A97").Value = "P3I"
$ws.Range("A98").Value = "P8I"
$ws.Range("A99").Value = "P14I"
$ws.Range("A100").Value = "P65I"
$ws.Range("A101").Value = "P66I"

$newRows = 97..101
foreach ($r in $newRows) {
    $ws.Range("D$r").Value = "imports"
    $ws.Range("E$r").Value = "PNW"
    $ws.Range("G$r").Value = 0
    $ws.Range("H$r").Value = 0
    $ws.Range("I$r").Value = 0
    $ws.Range("J$r").Value = 0
    $ws.Range("L$r").Value = 1
    $ws.Range("M$r").Value = 1
    $ws.Range("N$r").Value = 0
    $ws.Range("O$r").Value = 0
    $ws.Range("P$r").Value = 0
}

# --- view: scroll/selection follow the newly appended rows ---
$excel.Goto($ws.Range("A70"))
$ws.Range("A100").Select()
